$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$vals = @(17.54375623516665, 17.10651787051563, 16.83453691543235, 16.7229779492127, 16.7044147155641, 16.83303510957246, 17.3938190704251, 18.45904510052651, 19.21247584623408, 19.54735379748817, 19.67291787956133, 19.64593261599057, 19.55770965278964, 19.50350486375791, 19.19042241264254, 18.99625698116226, 18.88384453658146, 18.84566132295525, 19.0170030247689, 19.5836575966775, 19.94668387212475, 19.75363522471872, 19.0076261719345, 18.17544452833878)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $vals[$i]
}

# Column C
$vals = @(8.291723072200645, 8.014301920569171, 7.837640444863523, 7.764129406292616, 7.751833340146128, 7.836655105894542, 8.197420566601179, 8.851988505167123, 9.297495383818569, 9.491936538897658, 9.564345893284958, 9.54880612812255, 9.497918382035483, 9.466588048981494, 9.284619301170908, 9.170852144612457, 9.104644368437414, 9.08209624549171, 9.183043041216495, 9.512898777918858, 9.721345292172675, 9.61075776649311, 9.177534030462875, 8.680924003903844)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $vals[$i]
}

# Column D
$vals = @(8.862803453614745, 8.863698768900733, 8.865240734378737, 8.866119440578551, 8.866280492487242, 8.865251570219449, 8.862906736196305, 8.866143301718258, 8.873244079532014, 8.877487366762686, 8.879238684286838, 8.878855100986437, 8.877628560322202, 8.876896046967687, 8.872987055353709, 8.870847658583044, 8.869712568968554, 8.869344673862093, 8.871065532420722, 8.877984913928042, 8.883348645473097, 8.880409337598456, 8.870966736045274, 8.864434223788034)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $vals[$i]
}

# Column F
$vals = @(33.39431947920421, 33.44400992790391, 33.48257213437692, 33.50030780296821, 33.50337475975314, 33.48280314499744, 33.40977964098771, 33.33059847659814, 33.3115836277232, 33.31144836977278, 33.31262101928299, 33.31231405424857, 33.31152032240269, 33.31119349163978, 33.31176377906004, 33.31429434893165, 33.31655139639001, 33.31745325938106, 33.31394201108129, 33.31172025387738, 33.31740097645184, 33.31371682445008, 33.31409880421139, 33.34514866255619)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $vals[$i]
}

# Column G
$vals = @(3.656802260834634, 3.658894374284435, 3.660248150310283, 3.660817283143715, 3.660912843254176, 3.660255755075267, 3.657509290018005, 3.652670135678304, 3.649444576822486, 3.648048051079962, 3.647529347429825, 3.647640609787384, 3.64800517427771, 3.648229798340964, 3.649537263382934, 3.650357447836098, 3.650835862979868, 3.650998992682438, 3.650269448183695, 3.647897818313274, 3.646406845455346, 3.647197221479295, 3.650309211391142, 3.653921093987195)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $vals[$i]
}

# Column I
$vals = @(24.42223332966396, 24.52210522545123, 24.58830564304525, 24.61650810233125, 24.6212650569529, 24.58868103230427, 24.45565567248037, 24.23357173496857, 24.09414046586419, 24.03588353984268, 24.01456854631868, 24.0191259163212, 24.03411498700578, 24.0433933908448, 24.09805192281066, 24.13290904869643, 24.15344454692062, 24.16048102576017, 24.12914807189359, 24.0296920807922, 23.96903895657571, 24.00101223078909, 24.13084686564947, 24.289490232622)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $vals[$i]
}

# Column J
$vals = @(10.07989191966975, 10.10326352616297, 10.1184161802894, 10.12479329311239, 10.1258644404223, 10.1185013645398, 10.08778422324424, 10.03389113311867, 9.998129927601674, 9.982686788022924, 9.976956958864958, 9.978185732404437, 9.982213026158938, 9.984695233434426, 9.999155731465327, 10.0082377079067, 10.01353908260974, 10.01534738945761, 10.00726288158563, 9.981026909327065, 9.964568661628812, 9.973289891276876, 10.00770335125834, 10.04779498543006)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $vals[$i]
}

# Column L
$vals = @(12.02409148937367, 12.00897689313164, 12.0012944285417, 11.9985682054483, 11.99814001814785, 12.00125602085402, 12.01854995326313, 12.06501968091948, 12.10664039028352, 12.12715724480126, 12.13515024495266, 12.1334189251214, 12.12781037368728, 12.12440398654311, 12.10533109205658, 12.09403326723209, 12.08768434383165, 12.08556046541774, 12.09522051753156, 12.12945170238259, 12.15312574444198, 12.14037269607731, 12.09468330552195, 12.05112139130335)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $vals[$i]
}

# Column N
$vals = @(17.9419361833965, 17.98878498711294, 18.01933532078366, 18.03223453731328, 18.03440362918052, 18.01950746213387, 17.95771969346982, 17.85067856394008, 17.78059327482861, 17.7505570065319, 17.73944769910748, 17.74182852155959, 17.74963773433591, 17.75445556499873, 17.78259330256312, 17.80032719153402, 17.81070103662781, 17.81424330987586, 17.79842140791508, 17.74733679847161, 17.71549307572117, 17.7323476959089, 17.79928245693549, 17.87812944122966)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $vals[$i]
}

# Column O
$vals = @(25.26905837729013, 25.32733137593706, 25.36888763290723, 25.3872714275546, 25.39041147507594, 25.36912969962978, 25.28795015533492, 25.17472381366107, 25.11972369710967, 25.10085150861879, 25.09459077424144, 25.0958997180338, 25.1003186739388, 25.10314080775953, 25.12108090241121, 25.13366248673137, 25.14147764793415, 25.144223037315, 25.13226326228289, 25.0989966685881, 25.08241839155402, 25.09079361547469, 25.13289403937316, 25.20041509186001)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $vals[$i]
}
